$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; everything currently on row 1 (the
# "Archetype ID" / IM-name header row) shifts down to row 2, and all data
# rows shift down by one as well.
$ws.Rows.Item(1).Insert()

# New row 1 becomes a title row: "Archetype ID " (re-asserted, it will be
# merged vertically with the cell below it) in A1, and a new title
# "Percent Sufficiency- MVN GAM Model" spanning B1:I1.
$ws.Range("A1").Value() = "Archetype ID "
$ws.Range("B1").Value() = "Percent Sufficiency- MVN GAM Model"

# Give the new header row the same look as the existing header row
# (bold white text on a dark fill, centered) by copying that row's
# format onto it.
$ws.Range("A2:I2").Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Merge the title cells: the Archetype ID label spans both header rows,
# and the new title spans the remaining columns.
$ws.Range("A1:A2").Merge()
$ws.Range("B1:I1").Merge()

# Match the recorded selection / dimension state after the edit.
$ws.Range("J12").Select()
